$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 249.4252133562813
$ws.Range("B4").Value = 230.2965364763066
$ws.Range("B5").Value = 26.3923650749424
$ws.Range("B6").Value = 0
